# Rename the sheet "SQL Results" -> "Elenco"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Elenco"

# Force text formatting on the CODISTAT / IDREGIONE columns so that
# values with leading zeros (e.g. "072004", "03", "08") are kept as text
# instead of being coerced to numbers.
$ws.Range("A17:A19").NumberFormat = "@"
$ws.Range("E17:E19").NumberFormat = "@"

# Append three new rows of data (rows 17-19)
$ws.Range("A17").Value = "072004"
$ws.Range("B17").Value = "A225"
$ws.Range("C17").Value = "ALTAMURA"
$ws.Range("D17").Value = 72
$ws.Range("E17").Value = "16"
$ws.Range("F17").Value = "BA"
$ws.Range("G17").Value = 43025.5

$ws.Range("F18").Value = "MB"
$ws.Range("C18").Value = "BIASSONO"
$ws.Range("B18").Value = "A849"
$ws.Range("A18").Value = "108009"
$ws.Range("D18").Value = 108
$ws.Range("E18").Value = "03"
$ws.Range("G18").Value = 43027

$ws.Range("A19").Value = "037037"
$ws.Range("B19").Value = "F083"
$ws.Range("C19").Value = "MEDICINA"
$ws.Range("D19").Value = 37
$ws.Range("E19").Value = "08"
$ws.Range("F19").Value = "BO"
$ws.Range("G19").Value = 43028

# Copy formatting/number-format/border from the row above (row 16) into the new rows
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G19").PasteSpecial(-4122)  # xlPasteFormats

# Set the active selection to match the target state
$ws.Range("F27").Select()
